$d = $word.ActiveDocument

function Find-ParagraphByText($SearchText) {
    $rng = $d.Content
    $rng.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $rng.Paragraphs(1)
    $para
}

function Set-ParagraphBodyXml($Paragraph, $InnerXml) {
    $full = $Paragraph.Range
    # Exclude the trailing paragraph mark from the target range so only the
    # run content is swapped out; the paragraph's own pPr / identity
    # (paraId, textId, rsidR, ...) is left untouched by InsertXML.
    $target = $d.Range($full.Start, $full.End - 1)
    $target.InsertXML('<w:p>' + $InnerXml + '</w:p>')
}

# --- "Upload book" -> "Upload book" + " (done)" as a separate run ---
$uploadPara = Find-ParagraphByText "Upload book"
$uploadXml = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Upload book</w:t></w:r>'
$uploadXml = $uploadXml + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (done)</w:t></w:r>'
Set-ParagraphBodyXml $uploadPara $uploadXml

# --- "Read book" -> "Read book" + " (done)" as a separate run ---
$readPara = Find-ParagraphByText "Read book"
$readXml = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Read book</w:t></w:r>'
$readXml = $readXml + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (done)</w:t></w:r>'
Set-ParagraphBodyXml $readPara $readXml

# --- "Get google books from api" -> split "api" into its own run wrapped
#     in spellStart/spellEnd proofErr markers (as Word's proofer would do) ---
$apiPara = Find-ParagraphByText "Get google books from api"
$apiXml = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Get google books from </w:t></w:r>'
$apiXml = $apiXml + '<w:proofErr w:type="spellStart"/>'
$apiXml = $apiXml + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>api</w:t></w:r>'
$apiXml = $apiXml + '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphBodyXml $apiPara $apiXml
